# Update conversion "results" values (column G) for the Units reference list.
# Each row is keyed by its cell address; values are cast from string literals
# to doubles so Excel stores them canonically (as integers or in scientific
# notation, matching however the underlying magnitude requires).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Units")

$updates = @{
    "G6"   = [double]"1872605741056"
    "G13"  = [double]"6.743546251046112E+35"
    "G14"  = [double]"6.743546251046112E+35"
    "G23"  = [double]"5.872702763720583E+21"
    "G25"  = [double]"19997367730176"
    "G26"  = [double]"1.033985265379328E+20"
    "G29"  = [double]"1.033985265379328E+20"
    "G45"  = [double]"951232176848896"
    "G47"  = [double]"9.07588905277554E+20"
    "G50"  = [double]"9.07588905277554E+20"
    "G51"  = [double]"951232176848896"
    "G66"  = [double]"5.872702763720583E+21"
    "G67"  = [double]"5.872702763720583E+21"
    "G76"  = [double]"6.743546251046112E+35"
    "G99"  = [double]"9.07588905277554E+20"
    "G100" = [double]"14482629722112"
    "G104" = [double]"6.743546251046112E+35"
    "G105" = [double]"19997367730176"
    "G113" = [double]"6.743546251046112E+35"
    "G114" = [double]"14482629722112"
    "G117" = [double]"6.743546251046112E+35"
    "G118" = [double]"6.743546251046112E+35"
    "G120" = [double]"6.743546251046112E+35"
    "G123" = [double]"1.033985265379328E+20"
    "G127" = [double]"1872605741056"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
